# Insert a new weekly price record before the existing row 99 ("Fruta /
# hortaliza, semanal"). This shifts all subsequent records down by one row
# (old row 99 -> new row 100, ..., old row 217 -> new row 218) and the sheet
# dimension grows from A1:R217 to A1:R218.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 99, pushing row 99..217 down to
# 100..218.
$ws.Rows.Item(99).Insert()

# Make the new D99 a date-number cell formatted like the other "Fecha"
# cells in column D before writing its value.
$ws.Range("D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the newly inserted row 99 with the new data record.
$ws.Range("A99").Value = 4
$ws.Range("B99").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C99").Value = "Los Lagos"
$ws.Range("D99").Value = 44601
$ws.Range("E99").Value = 10
$ws.Range("F99").Value = 100112003
$ws.Range("G99").Value = "Ajo"
$ws.Range("H99").Value = "Chino"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 15
$ws.Range("K99").Value = 21000
$ws.Range("L99").Value = 21000
$ws.Range("M99").Value = 21000
$ws.Range("N99").Value = "$/caja 10 kilos"
$ws.Range("O99").Value = "China"
$ws.Range("P99").Value = 2100
$ws.Range("Q99").Value = 10
$ws.Range("R99").Value = "Hortaliza"
